$wb = $excel.ActiveWorkbook

# --- Sheet ALC: market-price refresh ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1004.3571
$ws.Range("I9").Value = 1046.4546
$ws.Range("K9").Value = 1046.4546
$ws.Range("M9").Value = -877.4546

$ws.Range("H100").Value = 3578.9092
$ws.Range("I100").Value = 2596.5557
$ws.Range("K100").Value = 2596.5557
$ws.Range("M100").Value = -2055.5557


# --- Sheet ARM: market-price refresh ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 488.8
$ws.Range("I4").Value = 450
$ws.Range("J4").Value = 514.6667
$ws.Range("K4").Value = 450
$ws.Range("L4").Value = 514.6667
$ws.Range("M4").Value = -334
$ws.Range("N4").Value = -746.6667

$ws.Range("H61").Value = 10571.556
$ws.Range("I61").Value = 10878
$ws.Range("J61").Value = 9499
$ws.Range("K61").Value = 10878
$ws.Range("L61").Value = 9499
$ws.Range("M61").Value = -10666
$ws.Range("N61").Value = -9923

$ws.Range("H122").Value = 2096.8333
$ws.Range("I122").Value = 2346.6
$ws.Range("J122").Value = 848
$ws.Range("K122").Value = 7039.799999999999
$ws.Range("L122").Value = 2544
$ws.Range("M122").Value = -4589.799999999999
$ws.Range("N122").Value = -7444

$ws.Range("H136").Value = 10571.556
$ws.Range("I136").Value = 10878
$ws.Range("J136").Value = 9499
$ws.Range("K136").Value = 32634
$ws.Range("L136").Value = 28497
$ws.Range("M136").Value = -30084
$ws.Range("N136").Value = -33597


# --- Sheet CRP: market-price refresh ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 186.33333
$ws.Range("I7").Value = 143.66667
$ws.Range("J7").Value = 399.66666
$ws.Range("K7").Value = 143.66667
$ws.Range("L7").Value = 399.66666
$ws.Range("M7").Value = -30.66667000000001
$ws.Range("N7").Value = -625.66666

$ws.Range("H16").Value = 4233
$ws.Range("I16").Value = 1350
$ws.Range("K16").Value = 1350
$ws.Range("M16").Value = -1063

$ws.Range("H31").Value = 5923.625
$ws.Range("I31").Value = 3314.8333
$ws.Range("J31").Value = 13750
$ws.Range("K31").Value = 3314.8333
$ws.Range("L31").Value = 13750
$ws.Range("M31").Value = -3019.8333
$ws.Range("N31").Value = -14340

$ws.Range("H34").Value = 5923.625
$ws.Range("I34").Value = 3314.8333
$ws.Range("J34").Value = 13750
$ws.Range("K34").Value = 3314.8333
$ws.Range("L34").Value = 13750
$ws.Range("M34").Value = -3112.8333
$ws.Range("N34").Value = -14154

$ws.Range("H58").Value = 2835.4707
$ws.Range("I58").Value = 1886
$ws.Range("K58").Value = 1886
$ws.Range("M58").Value = -1683

$ws.Range("H113").Value = 4233
$ws.Range("I113").Value = 1350
$ws.Range("K113").Value = 1350
$ws.Range("M113").Value = 820

$ws.Range("H136").Value = 2835.4707
$ws.Range("I136").Value = 1886
$ws.Range("K136").Value = 5658
$ws.Range("M136").Value = -3108


# --- Sheet CUL: market-price refresh ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3281509.5
$ws.Range("I4").Value = 1150147.5
$ws.Range("J4").Value = 15004000
$ws.Range("K4").Value = 3450442.5
$ws.Range("L4").Value = 45012000
$ws.Range("M4").Value = -3450330.5
$ws.Range("N4").Value = -45012224

$ws.Range("H109").Value = 2310
$ws.Range("I109").Value = 2750
$ws.Range("J109").Value = 2242.3076
$ws.Range("K109").Value = 8250
$ws.Range("L109").Value = 6726.9228
$ws.Range("M109").Value = -7210
$ws.Range("N109").Value = -8806.9228

$ws.Range("H124").Value = 11000
$ws.Range("I124").Value = 11000
$ws.Range("K124").Value = 33000
$ws.Range("M124").Value = -28090

$ws.Range("H126").Value = 9722.5
$ws.Range("J126").Value = 14495
$ws.Range("L126").Value = 43485
$ws.Range("N126").Value = -53365

$ws.Range("H131").Value = 4658.6
$ws.Range("I131").Value = 4658.6
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 13975.8
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = -8935.800000000001
$ws.Range("N131").ClearContents()


# --- Sheet GSM: market-price refresh ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2568.6667
$ws.Range("I80").Value = 3177.25
$ws.Range("K80").Value = 3177.25
$ws.Range("M80").Value = -2179.25

$ws.Range("H83").Value = 2568.6667
$ws.Range("I83").Value = 3177.25
$ws.Range("K83").Value = 15886.25
$ws.Range("M83").Value = -10894.25

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()


# --- Sheet LTW: market-price refresh ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

$ws.Range("H56").Value = 5400
$ws.Range("I56").Value = 5400
$ws.Range("K56").Value = 5400
$ws.Range("M56").Value = -4709

$ws.Range("H132").Value = 2205.6538
$ws.Range("I132").Value = 2206.125
$ws.Range("J132").Value = 2200
$ws.Range("K132").Value = 6618.375
$ws.Range("L132").Value = 6600
$ws.Range("M132").Value = -4088.375
$ws.Range("N132").Value = -11660


# --- Sheet WVR: market-price refresh ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()

$ws.Range("H61").Value = 144508.5
$ws.Range("I61").Value = 211262.75
$ws.Range("J61").Value = 11000
$ws.Range("K61").Value = 211262.75
$ws.Range("L61").Value = 11000
$ws.Range("M61").Value = -210970.75
$ws.Range("N61").Value = -11584

$ws.Range("H122").Value = 6541.4
$ws.Range("I122").Value = 6541.4
$ws.Range("K122").Value = 19624.2
$ws.Range("M122").Value = -17174.2

$ws.Range("H126").Value = 1115.75
$ws.Range("I126").Value = 989.9091
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 2969.7273
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -499.7273
$ws.Range("N126").Value = -12440

$ws.Range("H132").Value = 2709.5144
$ws.Range("I132").Value = 2504.0303
$ws.Range("K132").Value = 7512.090899999999
$ws.Range("M132").Value = -4982.090899999999

$ws.Range("H136").Value = 13047
$ws.Range("I136").Value = 10491.818
$ws.Range("J136").Value = 27100.5
$ws.Range("K136").Value = 31475.454
$ws.Range("L136").Value = 81301.5
$ws.Range("M136").Value = -28925.454
$ws.Range("N136").Value = -86401.5

